$d = $word.ActiveDocument
$nbh = [char]30

# ---------------------------------------------------------------------------
# Edit 1: "Important: The add-in makes use of the system clipboard."
#   -> "Important: PpspliT releases 1.27 and older makes use of the system
#       clipboard."
# ---------------------------------------------------------------------------
$find1 = $d.Content.Find
$needle1 = "The add" + $nbh + "in makes use of the system clipboard"
$ok1 = $find1.Execute($needle1, $true, $false, $false, $false, $false, $true, 1, $false, `
    "PpspliT releases 1.27 and older makes use of the system clipboard", 2)
Write-Output "edit1 (clipboard sentence): $ok1"

# ---------------------------------------------------------------------------
# Edit 2: insert a new highlighted paragraph right after the paragraph that
# ends in "...in operation." (and before the following empty paragraph),
# reusing the same paragraph style / border / shading as the surrounding
# "Important" callouts.
# ---------------------------------------------------------------------------
$find2 = $d.Content.Find
$ok2 = $find2.Execute("also interfere with the add" + $nbh + "in operation.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Output "edit2 (find anchor paragraph): $ok2"
$srcPara = $find2.Parent.Paragraphs(1)
$srcPara.Range.InsertParagraphAfter()
$newPara = $srcPara.Next()
$newRange = $newPara.Range
$newRange.Collapse(1)
$newRange.InsertAfter("Effective PpspliT release 2.0, this condition does not apply any longer: therefore, the system clipboard")
$newRange.Collapse(0)
$newRange.InsertAfter(" is left untouched and")
$newRange.Collapse(0)
$newRange.InsertAfter(" can be safely used even while a slide deck is being split.")
Write-Output "edit2 (new paragraph inserted)"

# ---------------------------------------------------------------------------
# Edit 3: "...unconfigured". M|ost notably|, this likely happens |when a
# pre-installed OEM PowerPoint release is being used|."
#   -> "...unconfigured". Most notably, this likely happens when a
#       pre-installed OEM PowerPoint release is being used."
# (merge the surrounding runs; the text itself does not change)
# ---------------------------------------------------------------------------
$find3 = $d.Content.Find
$ok3a = $find3.Execute("unconfigured", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Output "edit3 (find unconfigured): $ok3a"
$find3.Parent.Collapse(0)
$needle3a = " Most notably, this likely happens when a pre"
$ok3b = $find3.Execute($needle3a, $true, $false, $false, $false, $false, $true, 1, $false, $needle3a, 2)
Write-Output "edit3 (merge before hyphen): $ok3b"

$find3c = $d.Content.Find
$needle3c = "installed OEM PowerPoint release is being used."
$ok3c = $find3c.Execute($needle3c, $true, $false, $false, $false, $false, $true, 1, $false, $needle3c, 2)
Write-Output "edit3 (merge after hyphen): $ok3c"

# ---------------------------------------------------------------------------
# Edit 4: merge the split hyperlink display text into a single run. Editing
# the Hyperlink's TextToDisplay keeps the hyperlink field & character style
# intact (unlike a plain Find/Replace which would reset run formatting).
# ---------------------------------------------------------------------------
$mergedAny = $false
for ($i = 1; $i -le $d.Hyperlinks.Count; $i++) {
    $h = $d.Hyperlinks($i)
    if ($h.TextToDisplay -eq "https://support.microsoft.com/en-us/office/add-or-load-a-powerpoint-add-in-3de8bbc2-2481-457a-8841-7334cd5b455f") {
        $h.TextToDisplay = "https://support.microsoft.com/en-us/office/add-or-load-a-powerpoint-add-in-3de8bbc2-2481-457a-8841-7334cd5b455f"
        $mergedAny = $true
    }
}
Write-Output "edit4 (hyperlink text merged): $mergedAny"

# ---------------------------------------------------------------------------
# Edit 5: "For PowerPoint for Windows, c|lick the " -> "For PowerPoint for
# Windows, click the "
# ---------------------------------------------------------------------------
$find5 = $d.Content.Find
$needle5 = "For PowerPoint for Windows, click the"
$ok5 = $find5.Execute($needle5, $true, $false, $false, $false, $false, $true, 1, $false, $needle5, 2)
Write-Output "edit5: $ok5"

# ---------------------------------------------------------------------------
# Edit 6: " dialog box, browse for the |previously saved " -> " dialog box,
# browse for the previously saved "
# ---------------------------------------------------------------------------
$find6 = $d.Content.Find
$needle6 = "dialog box, browse for the previously saved"
$ok6 = $find6.Execute($needle6, $true, $false, $false, $false, $false, $true, 1, $false, $needle6, 2)
Write-Output "edit6: $ok6"

# ---------------------------------------------------------------------------
# Edit 7: " file and |then click " -> " file and then click "
# ---------------------------------------------------------------------------
$find7 = $d.Content.Find
$needle7 = "file and then click"
$ok7 = $find7.Execute($needle7, $true, $false, $false, $false, $false, $true, 1, $false, $needle7, 2)
Write-Output "edit7: $ok7"

# ---------------------------------------------------------------------------
# Edit 8: "A security notice |might appear|. |In this case just| click |on "
#   -> "A security notice might appear. In this case just click on "
# ---------------------------------------------------------------------------
$find8 = $d.Content.Find
$needle8 = "A security notice might appear. In this case just click on"
$ok8 = $find8.Execute($needle8, $true, $false, $false, $false, $false, $true, 1, $false, $needle8, 2)
Write-Output "edit8: $ok8"

# ---------------------------------------------------------------------------
# Edit 9: "%APPDATA%|\Microsoft\AddIns\PPspliT" -> "%APPDATA%\Microsoft\AddIns\PPspliT"
# ---------------------------------------------------------------------------
$find9 = $d.Content.Find
$needle9 = "%APPDATA%\Microsoft\AddIns\PPspliT"
$ok9 = $find9.Execute($needle9, $true, $false, $false, $false, $false, $true, 1, $false, $needle9, 2)
Write-Output "edit9: $ok9"

Write-Output "all edits applied"
